$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C14").Value2 = 44636
$ws.Range("D14").Value2 = 0.4375
